# Executed Login and Sign up script by testng.xml file
#
# Updates the sample sign-up data on the "SignUp" sheet: swaps the old
# test user rows (Rahul/deppak, Raj/rajnish, Geta/getta) for a fresh set
# (Sehat/sehat, abhishek/abhishek, narayn/narayan), nudges column B's
# width to fit the new values, and leaves the selection on C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignUp")

# --- Row 2: Sehat / sehat@asite.com ---
$ws.Range("A2").Value = "Sehat"
$ws.Range("B2").Value = "sehat@asite.com"

# --- Row 3: abhishek / abhishek@asite.com ---
$ws.Range("A3").Value = "abhishek"
$ws.Range("B3").Value = "abhishek@asite.com"

# --- Row 4: narayn / narayan@asite.com ---
$ws.Range("A4").Value = "narayn"
$ws.Range("B4").Value = "narayan@asite.com"

# Column B picked up a (no-op) fill touch in the source session -- mirror
# it so the cells carry the same applied-fill formatting flag.
$ws.Range("B2:B4").Interior.Pattern = -4142
$ws.Range("B2:B4").Interior.ColorIndex = -4142

# Column B widened slightly to fit the new e-mail addresses.
$ws.Columns.Item(2).ColumnWidth = 17.29

# Final selection left on C8.
$ws.Range("C8").Select()
